$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 23.42000000000022
$ws.Range("H2").Value = 0.0000000000000002203916674193859
$ws.Range("K2").Value = 41.62679280669325
$ws.Range("L2").Value = "[34.25339512596389, 49.0001904874226]"
$ws.Range("O2").Value = 1.415131825941349
$ws.Range("P2").Value = "[1.2264475824825025, 1.603816069400195]"
$ws.Range("S2").Value = 64.25437799752753
$ws.Range("T2").Value = "[59.79663405302932, 68.71212194202573]"
$ws.Range("W2").Value = 18.1452252252254
$ws.Range("X2").Value = 17.44192192192209
$ws.Range("Y2").Value = 18.8485285285287

# Row 3 updates
$ws.Range("B3").Value = 1
$ws.Range("E3").Value = 22.79000000000012
$ws.Range("H3").Value = 0.0000000000000002203916674193859
$ws.Range("K3").Value = 48.54592058412872
$ws.Range("L3").Value = "[40.981783575249324, 56.11005759300811]"
$ws.Range("O3").Value = -0.2012631930227693
$ws.Range("P3").Value = "[-0.3647895373537695, -0.037736848691769076]"
$ws.Range("Q3").Value = 0.01598976796917784
$ws.Range("R3").Value = 0.01598976796917784
$ws.Range("S3").Value = 63.34235063271495
$ws.Range("T3").Value = "[58.67971057150177, 68.00499069392814]"
$ws.Range("W3").Value = 0.7300100100100124
$ws.Range("X3").Value = 0.1368768768768761
$ws.Range("Y3").Value = 1.323143143143149
